$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.298.12"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.23%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.813.71"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +3.46%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.06%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'326.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.83%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9994"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.11%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4364"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.07%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3675"
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = "'44.72"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -1.08%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.07658"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +2.25%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'1.141"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.71%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.17%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'21.99"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.72%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = "'  +2.74%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'7.507"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +3.73%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'1.812.56"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +3.92%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'95.63"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +8.99%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.01%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.06530"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +5.07%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.9996"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.11%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'17.38"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.93%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'6.252"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.65%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'28.311.60"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.27%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'11.55"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.49%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.111"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -9.66%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'161.95"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +6.38%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'20.75"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.99%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'2.021.14"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +3.92%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.283"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -3.53%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'129.04"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.30%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'1.202"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -2.24%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'5.991"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.17%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.09163"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.06%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'3.486"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -5.21%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'12.99"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +1.92%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.35%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.2174"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +1.14%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'5.184"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.27%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.6589"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.60%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.06204"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.38%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.194"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -0.38%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'8.121"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.87%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'1.426"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.34%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.9993"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.12%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'13.87"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.61%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.6117"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +2.92%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'3.740"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.10%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'125.65"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.38%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.017"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.31%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'Cronos"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'0.07001"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.49%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'EOS"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'1.154"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.65%  "
$ws.Range('E51').Style = 'Normal'
